$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 89 (pushes the existing row 89 and all
# subsequent rows down by one, growing the table from 141 to 142 data rows).
$ws.Rows.Item(89).EntireRow.Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A89").Value = 6
$ws.Range("B89").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C89").Value = "Metropolitana"
$ws.Range("D89").Value = 44572
$ws.Range("E89").Value = 13
$ws.Range("F89").Value = 100112029
$ws.Range("G89").Value = "Orégano"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 34
$ws.Range("K89").Value = 8500
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = 8721
$ws.Range("N89").Value = "$/docena de atados"
$ws.Range("O89").Value = "Región Metropolitana"
$ws.Range("P89").Value = 2907
$ws.Range("Q89").Value = 3
$ws.Range("R89").Value = "Hortaliza"
